$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''64.342.93'
$ws.Range("E2").Value = '''  -2.84%  '
$ws.Range("D3").Value = '''3.148.39'
$ws.Range("E3").Value = '''  -1.72%  '
$ws.Range("E4").Value = '''  +0.05%  '
$ws.Range("D5").Value = '''608.27'
$ws.Range("E5").Value = '''  +0.19%  '
$ws.Range("D6").Value = '''149.50'
$ws.Range("E6").Value = '''  -4.23%  '
$ws.Range("E7").Value = '''  +0.10%  '
$ws.Range("D8").Value = '''3.146.75'
$ws.Range("E8").Value = '''  -1.71%  '
$ws.Range("E9").Value = '''  -3.31%  '
$ws.Range("E10").Value = '''  -4.42%  '
$ws.Range("D11").Value = '''5.63'
$ws.Range("E11").Value = '''  -0.34%  '
$ws.Range("E12").Value = '''  -4.48%  '
$ws.Range("D14").Value = '''36.87'
$ws.Range("E14").Value = '''  -4.00%  '
$ws.Range("D15").Value = '''3.663.28'
$ws.Range("D16").Value = '''64.360.09'
$ws.Range("E16").Value = '''  -3.04%  '
$ws.Range("B17").Value = '''TRON'
$ws.Range("C17").Value = '''https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '''0.114'
$ws.Range("E17").Value = '''  +0.12%  '
$ws.Range("B18").Value = '''WrappedEther'
$ws.Range("C18").Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '''3.142.98'
$ws.Range("E18").Value = '''  -1.98%  '
$ws.Range("E19").Value = '''  -3.14%  '
$ws.Range("D20").Value = '''483.62'
$ws.Range("D21").Value = '''14.69'
$ws.Range("E21").Value = '''  -3.98%  '
$ws.Range("E22").Value = '''  -2.06%  '
$ws.Range("E23").Value = '''  -2.60%  '
$ws.Range("D24").Value = '''13.86'
$ws.Range("E24").Value = '''  -5.26%  '
$ws.Range("D25").Value = '''84.22'
$ws.Range("E25").Value = '''  -1.04%  '
$ws.Range("E26").Value = '''  +0.13%  '
$ws.Range("D27").Value = '''2.95'
$ws.Range("E27").Value = '''  -1.46%  '
$ws.Range("D28").Value = '''8.58'
$ws.Range("E28").Value = '''  -4.95%  '
$ws.Range("E29").Value = '''  -4.03%  '
$ws.Range("D30").Value = '''0.126'
$ws.Range("E30").Value = '''  -3.96%  '
$ws.Range("D31").Value = '''6.98'
$ws.Range("E31").Value = '''  +0.43%  '
$ws.Range("D32").Value = '''2.74'
$ws.Range("E32").Value = '''  -6.04%  '
$ws.Range("E33").Value = '''  -0.19%  '
$ws.Range("D34").Value = '''26.80'
$ws.Range("E34").Value = '''  -4.88%  '
$ws.Range("E35").Value = '''  -5.30%  '
$ws.Range("E36").Value = '''  -4.60%  '
$ws.Range("D37").Value = '''54.44'
$ws.Range("E37").Value = '''  -1.66%  '
$ws.Range("E38").Value = '''  +6.94%  '
$ws.Range("E39").Value = '''  -1.75%  '
$ws.Range("D40").Value = '''453.55'
$ws.Range("E40").Value = '''  -9.23%  '
$ws.Range("D41").Value = '''0.0402'
$ws.Range("E41").Value = '''  -3.91%  '
$ws.Range("E42").Value = '''  -5.52%  '
$ws.Range("D43").Value = '''8.50'
$ws.Range("E43").Value = '''  -2.43%  '
$ws.Range("D44").Value = '''2.894.68'
$ws.Range("E44").Value = '''  -0.52%  '
$ws.Range("E45").Value = '''  -7.50%  '
$ws.Range("D46").Value = '''2.33'
$ws.Range("E46").Value = '''  -4.09%  '
$ws.Range("D47").Value = '''26.80'
$ws.Range("E47").Value = '''  -4.74%  '
$ws.Range("D48").Value = '''0.999'
$ws.Range("E48").Value = '''  +0.04%  '
$ws.Range("B49").Value = '''ThetaToken'
$ws.Range("C49").Value = '''https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").Value = '''2.34'
$ws.Range("E49").Value = '''  -2.55%  '
$ws.Range("B50").Value = '''Stellar'
$ws.Range("C50").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '''0.115'
$ws.Range("E50").Value = '''  -1.07%  '
$ws.Range("D51").Value = '''120.22'
$ws.Range("E51").Value = '''  -1.58%  '
